$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, border, centered) from H1 to the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for columns I and J
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 5
